# Apply cryptos-list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '95.120.56'
$ws.Range("E2").Value = '  -1.09%  '

# Row 3
$ws.Range("D3").Value = '3.573.77'
$ws.Range("E3").Value = '  -1.60%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.33'
$ws.Range("E5").Value = '  -1.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '656.57'
$ws.Range("E6").Value = '  +2.55%  '

# Row 7
$ws.Range("E7").Value = '  -0.68%  '

# Row 8
$ws.Range("E8").Value = '  -0.14%  '

# Row 9
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
$ws.Range("E10").Value = '  +0.04%  '

# Row 11
$ws.Range("D11").Value = '3.572.60'
$ws.Range("E11").Value = '  -1.65%  '

# Row 12
$ws.Range("E12").Value = '  +1.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.46'
$ws.Range("E13").Value = '  -1.95%  '

# Row 14
$ws.Range("E14").Value = '  +1.80%  '

# Row 15
$ws.Range("D15").Value = '4.236.42'
$ws.Range("E15").Value = '  -2.12%  '

# Row 16
$ws.Range("D16").Value = '95.017.09'
$ws.Range("E16").Value = '  -1.18%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000254'
$ws.Range("E17").Value = '  -0.41%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.566.81'
$ws.Range("E18").Value = '  -2.10%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.06'
$ws.Range("E19").Value = '  +0.58%  '

# Row 20
$ws.Range("E20").Value = '  -5.15%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.86'
$ws.Range("E21").Value = '  -2.48%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.46'
$ws.Range("E22").Value = '  +0.82%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.76'
$ws.Range("E23").Value = '  -1.54%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.480'
$ws.Range("E24").Value = '  -3.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.97'
$ws.Range("E25").Value = '  +4.25%  '

# Row 26
$ws.Range("E26").Value = '  -0.08%  '

# Row 27
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.79'
$ws.Range("E27").Value = '  +2.48%  '

# Row 28
$ws.Range("B28").Value = 'Litecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '90.01'
$ws.Range("E28").Value = '  -7.48%  '

# Row 29
$ws.Range("D29").Value = '3.763.12'
$ws.Range("E29").Value = '  -1.80%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.04'
$ws.Range("E30").Value = '  -2.32%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.144'
$ws.Range("E31").Value = '  +1.13%  '

# Row 32
$ws.Range("E32").Value = '  +0.11%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.48%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.177'
$ws.Range("E35").Value = '  -2.17%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.04'
$ws.Range("E36").Value = '  +2.86%  '

# Row 37
$ws.Range("E37").Value = '  +17.69%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '605.37'
$ws.Range("E38").Value = '  +6.27%  '

# Row 39
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.561'
$ws.Range("E39").Value = '  -2.16%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.57'
$ws.Range("E40").Value = '  +9.22%  '

# Row 41
$ws.Range("E41").Value = '  -0.18%  '

# Row 42
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.909'
$ws.Range("E43").Value = '  -3.12%  '

# Row 44
$ws.Range("E44").Value = '  +6.72%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '35.12'
$ws.Range("E45").Value = '  +25.43%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.75'
$ws.Range("E46").Value = '  +0.39%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.31'
$ws.Range("E47").Value = '  +4.33%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.41'
$ws.Range("E48").Value = '  -1.62%  '

# Row 49
$ws.Range("E49").Value = '  -3.18%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.51'
$ws.Range("E50").Value = '  +0.05%  '

# Row 51
$ws.Range("E51").Value = '  +0.28%  '
